$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Coin name / Link / Volume text cells (safe as plain text assignment)
$ws.Range("E2").Value = "  -4.48%  "
$ws.Range("E3").Value = "  -4.52%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("E6").Value = "  -2.69%  "
$ws.Range("E7").Value = "  -2.96%  "
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("E9").Value = "  -3.93%  "
$ws.Range("E10").Value = "  -2.82%  "
$ws.Range("E11").Value = "  -5.90%  "
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("E13").Value = "  -3.88%  "
$ws.Range("E14").Value = "  -5.13%  "
$ws.Range("E15").Value = "  -5.51%  "
$ws.Range("E16").Value = "  -5.05%  "
$ws.Range("E17").Value = "  -4.04%  "
$ws.Range("E18").Value = "  -7.97%  "
$ws.Range("E19").Value = "  -2.93%  "
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("E21").Value = "  -7.34%  "
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("E22").Value = "  -5.12%  "
$ws.Range("B23").Value = "BitDAO"
$ws.Range("C23").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("E23").Value = "  -7.25%  "
$ws.Range("E24").Value = "  -2.78%  "
$ws.Range("E25").Value = "  -4.43%  "
$ws.Range("E26").Value = "  -3.52%  "
$ws.Range("E27").Value = "  +3.59%  "
$ws.Range("E28").Value = "  -4.69%  "
$ws.Range("E29").Value = "  -5.02%  "
$ws.Range("E30").Value = "  -4.43%  "
$ws.Range("E31").Value = "  -5.13%  "
$ws.Range("E32").Value = "  -4.07%  "
$ws.Range("E33").Value = "  +5.55%  "
$ws.Range("E34").Value = "  -9.25%  "
$ws.Range("E35").Value = "  -6.85%  "
$ws.Range("E36").Value = "  -10.17%  "
$ws.Range("E37").Value = "  -3.21%  "
$ws.Range("E38").Value = "  -6.12%  "
$ws.Range("E39").Value = "  -4.55%  "
$ws.Range("E41").Value = "  -6.56%  "
$ws.Range("E42").Value = "  -7.85%  "
$ws.Range("E43").Value = "  -3.23%  "
$ws.Range("E44").Value = "  -6.85%  "
$ws.Range("E45").Value = "  -7.99%  "
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("E47").Value = "  -5.29%  "
$ws.Range("E48").Value = "  -3.08%  "
$ws.Range("E49").Value = "  -5.01%  "
$ws.Range("E50").Value = "  +4.78%  "
$ws.Range("E51").Value = "  -2.92%  "

# Update Price cells - force text storage so numeric-looking strings are not
# auto-converted into real numbers (must remain text, matching source data)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.398.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "290.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3673"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.32"
$ws.Range("D8").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07597"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.053"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.901"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.567.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "89.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06773"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.228"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5304"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "22.417.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.399"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.995"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "145.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.982"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "125.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.744.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.039"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.256"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.991"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.34"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08449"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02545"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2326"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06542"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.526"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.250"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6386"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.0000"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5998"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.782"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.138"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.230"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "123.39"
$ws.Range("D51").Style = "Normal"
